# Updated the page names for Subpopulation and LOT pages.
#
# The workbook uses short identifier strings (in column F/G/H) as
# page/section names. Rename the "Subpopulation" and "LOT" identifiers to
# the new "pop_filter" naming scheme, wherever they occur in the used
# range of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "sub_pop_section"           = "pop_filter1_section"
    "sub_pop_section1"          = "pop_filter1_section1"
    "sub_pop_section1_checkbox" = "pop_filter1_section1_checkbox"
    "sub_pop_section2"          = "pop_filter1_section2"
    "sub_pop_section2_checkbox" = "pop_filter1_section2_checkbox"
    "lot_section"                = "pop_filter2_section"
    "lot_section2"                = "pop_filter2_section2"
    "lot_section2_checkbox"      = "pop_filter2_section2_checkbox"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $map.ContainsKey([string]$val)) {
            $cell.Value = $map[[string]$val]
        }
    }
}

# Match the final selection recorded in the workbook after the edit.
$ws.Range("H10").Select()
